$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B and H hold numeric-looking text ("8", "38.0", ...). Force Text number
# format first so COM does not auto-coerce the assigned strings into numbers
# (matches the workbook author storing them as text throughout the table).
$ws.Range("B6:B13").NumberFormat = "@"
$ws.Range("H6:H13").NumberFormat = "@"
$ws.Range("B16").NumberFormat = "@"

# Row 6
$ws.Range("B6").Value = "8"
$ws.Range("C6").Value = "Marni  "
$ws.Range("D6").Value = "Shanika  "
$ws.Range("E6").Value = "-2.63,7.34"
$ws.Range("F6").Value = "Lady(mother): 0560804012"
$ws.Range("G6").Value = "7:00:00"
$ws.Range("H6").Value = "38.0"

# Row 7
$ws.Range("B7").Value = "9"
$ws.Range("C7").Value = "Letha  "
$ws.Range("D7").Value = "Stephenie  "
$ws.Range("E7").Value = "-9.76,7.61"
$ws.Range("F7").Value = "Sibyl(mother): 0567328221"
$ws.Range("G7").Value = "7:08:00"
$ws.Range("H7").Value = "30.0"

# Row 8
$ws.Range("B8").Value = "11"
$ws.Range("C8").Value = "Randolph  "
$ws.Range("D8").Value = "Bridgette  "
$ws.Range("E8").Value = "-6.35,3.21"
$ws.Range("F8").Value = "Lenny(father): 0505536740"
$ws.Range("G8").Value = "7:14:00"
$ws.Range("H8").Value = "24.0"

# Row 9
$ws.Range("B9").Value = "19"
$ws.Range("C9").Value = "Jeanine  "
$ws.Range("D9").Value = "Janee  "
$ws.Range("E9").Value = "-7.76,-1.4"
$ws.Range("F9").Value = "Teresa(mother): 0517627420"
$ws.Range("G9").Value = "7:20:00"
$ws.Range("H9").Value = "18.0"

# Row 10
$ws.Range("B10").Value = "7"
$ws.Range("C10").Value = "Wyatt  "
$ws.Range("D10").Value = "Willette  "
$ws.Range("E10").Value = "-4.33,0.85"
$ws.Range("F10").Value = "Antionette(father): 0557331799"
$ws.Range("G10").Value = "7:27:00"
$ws.Range("H10").Value = "11.0"

# Row 11
$ws.Range("B11").Value = "13"
$ws.Range("C11").Value = "Fay  "
$ws.Range("D11").Value = "Emilee  "
$ws.Range("E11").Value = "-4.89,2.74"
$ws.Range("F11").Value = "Sheri(mother): 0516797453"
$ws.Range("G11").Value = "7:29:00"
$ws.Range("H11").Value = "9.0"

# Row 12
$ws.Range("B12").Value = "15"
$ws.Range("C12").Value = "Nubia  "
$ws.Range("D12").Value = "Royce  "
$ws.Range("E12").Value = "-5.03,3.47"
$ws.Range("F12").Value = "Augustus(father): 0517389040"
$ws.Range("G12").Value = "7:30:00"
$ws.Range("H12").Value = "8.0"

# Row 13
$ws.Range("B13").Value = "0"
$ws.Range("C13").Value = "Trudie  "
$ws.Range("D13").Value = "Fleta  "
$ws.Range("E13").Value = "-3.01,3.2"
$ws.Range("F13").Value = "Anneliese(father): 0548973345"
$ws.Range("G13").Value = "7:33:00"
$ws.Range("H13").Value = "5.0"

# Summary rows: school pickup time and total time-in-minutes
$ws.Range("G14").Value = "7:38:00"
$ws.Range("B16").Value = "38.0"
